# Add "skip transcoding" support to the batch-ingest manifest.
#
# The existing sheet has a header row (row 2) with File/Offset/Label columns
# (E:G) and one data row (row 3) under them, plus a stray empty-but-styled
# cell at F4. This change duplicates the File/Offset/Label pair into a new
# H:J block, adds a new "Skip Transcoding" column (K) with a "yes" value on
# the data row, and drops the now-unused row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2): H2:J2 mirror E2:G2 (File / Offset / Label) ---
$ws.Range("H2").Value = "File"
$ws.Range("E2").Copy()
$ws.Range("H2").PasteSpecial(-4122)   # xlPasteFormats - reuse E2's style (s=1)

$ws.Range("I2").Value = "Offset"
$ws.Range("F2").Copy()
$ws.Range("I2").PasteSpecial(-4122)   # reuse F2's style (s=1)

$ws.Range("J2").Value = "Label"
$ws.Range("G2").Copy()
$ws.Range("J2").PasteSpecial(-4122)   # reuse G2's style (s=1)
$excel.CutCopyMode = 0

# --- Data row (row 3): H3:I3 mirror E3:F3 (file path / offset) ---
$ws.Range("H3").Value = "assets/sheephead_mountain.mov"   # same style as E3 (default)

$ws.Range("I3").Value = "00:00:00.500"
$ws.Range("F3").Copy()
$ws.Range("I3").PasteSpecial(-4122)   # reuse F3's quote-prefixed style (s=3)
$excel.CutCopyMode = 0

# --- New content: Label for the duplicated block, then the new Skip
#     Transcoding column. Order matters so new shared strings are appended
#     in the same sequence as the authored workbook (J3, then K2, then K3).
$ws.Range("J3").Value = "Unde aliquid"                # new label text, default style

$ws.Range("K2").Value = "Skip Transcoding"
$ws.Range("E2").Copy()
$ws.Range("K2").PasteSpecial(-4122)   # header style (s=1)
$excel.CutCopyMode = 0

$ws.Range("K3").Value = "yes"                         # new data value, default style

# --- Drop the old stray cell at F4 so row 4 disappears entirely ---
$ws.Range("F4").Clear()

# --- Size the three new columns to fit their content ---
$ws.Columns.Item(8).ColumnWidth = 27.998697916666668   # -> ~28.83 chars wide
$ws.Columns.Item(9).ColumnWidth = 11.166666666666666   # -> 12 chars wide
$ws.Columns.Item(10).ColumnWidth = 10.666666666666666  # -> 11.5 chars wide

# --- Match the saved selection state ---
$ws.Range("B11").Select() | Out-Null
